$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.88
$ws.Range("J2").Value = 3.6
$ws.Range("N2").Value = 4.5
$ws.Range("Q2").Value = 1.8
$ws.Range("S2").Value = 3
$ws.Range("W2").Value = 1.53
$ws.Range("X2").Value = 18.5
$ws.Range("Z2").Value = 20
$ws.Range("AK2").Value = 32
$ws.Range("F5").Value = 3.9
$ws.Range("G5").Value = 5.1
$ws.Range("I5").Value = 1.98
$ws.Range("V5").Value = 2.02
$ws.Range("G6").Value = 3.7
$ws.Range("H6").Value = 2.2
$ws.Range("I6").Value = 2.22
$ws.Range("L6").Value = 1.41
$ws.Range("P6").Value = 2.02
$ws.Range("Q6").Value = 1.97
$ws.Range("V6").Value = 1.81
$ws.Range("X6").Value = 14
$ws.Range("Y6").Value = 10
$ws.Range("AB6").Value = 14
$ws.Range("AF6").Value = 27
$ws.Range("AH6").Value = 17.5
$ws.Range("AJ6").Value = 65
$ws.Range("AK6").Value = 42
$ws.Range("T7").Value = 1.9
$ws.Range("F8").Value = 2.54
$ws.Range("G8").Value = 2.58
$ws.Range("H8").Value = 2.92
$ws.Range("I8").Value = 3
$ws.Range("O8").Value = 1.01
$ws.Range("V8").Value = 1.5
$ws.Range("W8").Value = 1.63
$ws.Range("X8").Value = 650
$ws.Range("F9").Value = 2.44
$ws.Range("G9").Value = 3.9
$ws.Range("H9").Value = 2.02
$ws.Range("J9").Value = 2.36
$ws.Range("K9").Value = 950
$ws.Range("N9").Value = 1.02
$ws.Range("O9").Value = 1.19
$ws.Range("P9").Value = 2.02
$ws.Range("Q9").Value = 1.19
$ws.Range("R9").Value = 1.41
$ws.Range("S9").Value = 1.19
$ws.Range("W9").Value = 1.35
$ws.Range("F10").Value = 2.24
$ws.Range("G10").Value = 2.74
$ws.Range("H10").Value = 2.92
$ws.Range("I10").Value = 4.3
$ws.Range("K10").Value = 3.7
$ws.Range("N10").Value = 1.81
$ws.Range("O10").Value = 1.3
$ws.Range("P10").Value = 1.67
$ws.Range("Q10").Value = 1.87
$ws.Range("R10").Value = 1.22
$ws.Range("S10").Value = 2.84
$ws.Range("T10").Value = 1.01
$ws.Range("U10").Value = 1.8
$ws.Range("V10").Value = 1.3
$ws.Range("F11").Value = 3.1
$ws.Range("G11").Value = 3.8
$ws.Range("I11").Value = 2.94
$ws.Range("L11").Value = 1.41
$ws.Range("V11").Value = 1.55
$ws.Range("F12").Value = 1.62
$ws.Range("H12").Value = 5.2
$ws.Range("I12").Value = 6.4
$ws.Range("N12").Value = 4.3
$ws.Range("W12").Value = 2.44
$ws.Range("AB12").Value = 12
$ws.Range("F13").Value = 3
$ws.Range("G13").Value = 3.2
$ws.Range("H13").Value = 2.4
$ws.Range("I13").Value = 2.62
$ws.Range("K13").Value = 3.75
$ws.Range("L13").Value = 1.37
$ws.Range("S13").Value = 2.86
$ws.Range("L14").Value = 1.28
$ws.Range("M14").Value = 1.04
$ws.Range("L15").Value = 1.28
$ws.Range("U15").Value = 2.1
$ws.Range("F16").Value = 2.5
$ws.Range("G16").Value = 2.6
$ws.Range("H16").Value = 2.88
$ws.Range("I16").Value = 3.1
$ws.Range("J16").Value = 3.5
$ws.Range("K16").Value = 3.7
$ws.Range("U16").Value = 2.14
$ws.Range("V16").Value = 1.47
$ws.Range("W16").Value = 1.62
$ws.Range("F18").Value = 2.1
$ws.Range("G18").Value = 2.12
$ws.Range("I18").Value = 3.85
$ws.Range("J18").Value = 3.8
$ws.Range("K18").Value = 3.85
$ws.Range("L18").Value = 1.33
$ws.Range("N18").Value = 5
$ws.Range("O18").Value = 1.24
$ws.Range("P18").Value = 2.3
$ws.Range("R18").Value = 1.53
$ws.Range("S18").Value = 2.82
$ws.Range("W18").Value = 1.9
$ws.Range("X18").Value = 19
$ws.Range("AH18").Value = 15.5
$ws.Range("AK18").Value = 19.5
$ws.Range("F19").Value = 5.8
$ws.Range("G19").Value = 6
$ws.Range("H19").Value = 1.63
$ws.Range("L19").Value = 1.29
$ws.Range("N19").Value = 6
$ws.Range("P19").Value = 2.66
$ws.Range("V19").Value = 2.56
$ws.Range("W19").Value = 1.2
$ws.Range("AN19").Value = 50
$ws.Range("F20").Value = 3.9
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = 1.98
$ws.Range("I20").Value = 1.99
$ws.Range("L20").Value = 1.3
$ws.Range("P20").Value = 2.5
$ws.Range("S20").Value = 2.62
$ws.Range("U20").Value = 2.56
$ws.Range("V20").Value = 2
$ws.Range("W20").Value = 1.33
$ws.Range("X20").Value = 22
$ws.Range("Y20").Value = 13
$ws.Range("Z20").Value = 14
$ws.Range("AB20").Value = 20
$ws.Range("AD20").Value = 10
$ws.Range("AI20").Value = 26
$ws.Range("AK20").Value = 40
$ws.Range("AN20").Value = 29
$ws.Range("AO20").Value = 9.4
$ws.Range("G21").Value = 2.04
$ws.Range("L21").Value = 1.01
$ws.Range("M21").Value = 1.01
$ws.Range("N21").Value = 1.58
$ws.Range("O21").Value = 1.48
$ws.Range("R21").Value = 1.16
$ws.Range("S21").Value = 4.3
$ws.Range("T21").Value = 1.83
$ws.Range("U21").Value = 1.53
$ws.Range("V21").Value = 1.23
$ws.Range("W21").Value = 1.96
$ws.Range("X21").Value = 13
$ws.Range("Y21").Value = 18.5
$ws.Range("Z21").Value = 50
$ws.Range("AA21").Value = 1000
$ws.Range("AB21").Value = 9.199999999999999
$ws.Range("AC21").Value = 10.5
$ws.Range("AD21").Value = 29
$ws.Range("AE21").Value = 1000
$ws.Range("AF21").Value = 15
$ws.Range("AG21").Value = 15
$ws.Range("AH21").Value = 34
$ws.Range("AI21").Value = 1000
$ws.Range("AJ21").Value = 34
$ws.Range("AK21").Value = 36
$ws.Range("AL21").Value = 75
$ws.Range("AM21").Value = 1000
$ws.Range("AN21").Value = 1000
$ws.Range("AO21").Value = 1000
$ws.Range("I22").Value = 2.44
$ws.Range("J22").Value = 3.1
$ws.Range("U22").Value = 1.78
$ws.Range("AA22").Value = 34
